$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6512.1113
$ws.Range("I19").Value = 9183.799999999999
$ws.Range("J19").Value = 3172.5
$ws.Range("K19").Value = 9183.799999999999
$ws.Range("L19").Value = 3172.5
$ws.Range("M19").Value = -9008.799999999999
$ws.Range("N19").Value = -3522.5

$ws.Range("H28").Value = 5211.8184
$ws.Range("I28").Value = 506.7143
$ws.Range("J28").Value = 13445.75
$ws.Range("K28").Value = 506.7143
$ws.Range("L28").Value = 13445.75
$ws.Range("M28").Value = -21.71429999999998
$ws.Range("N28").Value = -14415.75

$ws.Range("H98").Value = 8687.777
$ws.Range("I98").Value = 8983.846
$ws.Range("J98").Value = 990
$ws.Range("K98").Value = 8983.846
$ws.Range("L98").Value = 990
$ws.Range("M98").Value = -7485.846
$ws.Range("N98").Value = -3986

$ws.Range("H107").Value = 9487.556
$ws.Range("I107").Value = 485.42856
$ws.Range("J107").Value = 40995
$ws.Range("K107").Value = 485.42856
$ws.Range("L107").Value = 40995
$ws.Range("M107").Value = 1434.57144
$ws.Range("N107").Value = -44835

$ws.Range("H111").Value = 2016.7693
$ws.Range("I111").Value = 2398.4285
$ws.Range("J111").Value = 1571.5
$ws.Range("K111").Value = 7195.2855
$ws.Range("L111").Value = 4714.5
$ws.Range("M111").Value = -4128.2855
$ws.Range("N111").Value = -10848.5

$ws.Range("H113").Value = 3007.6924
$ws.Range("I113").Value = 2425
$ws.Range("J113").Value = 3266.6667
$ws.Range("K113").Value = 2425
$ws.Range("L113").Value = 3266.6667
$ws.Range("M113").Value = 829
$ws.Range("N113").Value = -9774.6667

$ws.Range("H116").Value = 2837.1875
$ws.Range("J116").Value = 3309.0908
$ws.Range("L116").Value = 3309.0908
$ws.Range("N116").Value = -10193.0908

$ws.Range("H122").Value = 8687.777
$ws.Range("I122").Value = 8983.846
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 26951.538
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -24501.538
$ws.Range("N122").Value = -7870

$ws.Range("H137").Value = 23258544
$ws.Range("I137").Value = 1189.0769
$ws.Range("J137").Value = 58828616
$ws.Range("K137").Value = 3567.2307
$ws.Range("L137").Value = 176485848
$ws.Range("M137").Value = -1017.2307
$ws.Range("N137").Value = -176490948

$ws.Range("H138").Value = 3427.4
$ws.Range("I138").Value = 1263.3103
$ws.Range("J138").Value = 4311.3237
$ws.Range("K138").Value = 3789.9309
$ws.Range("L138").Value = 12933.9711
$ws.Range("M138").Value = 1350.0691
$ws.Range("N138").Value = -23213.9711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 72384.94500000001
$ws.Range("I32").Value = 75558.06
$ws.Range("J32").Value = 47000
$ws.Range("K32").Value = 75558.06
$ws.Range("L32").Value = 47000
$ws.Range("M32").Value = -75271.06
$ws.Range("N32").Value = -47574

$ws.Range("H61").Value = 2179.8262
$ws.Range("I61").Value = 1979.7778
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 1979.7778
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -1767.7778
$ws.Range("N61").Value = -3324

$ws.Range("H74").Value = 1460.1864
$ws.Range("I74").Value = 1400.3864
$ws.Range("J74").Value = 1635.6
$ws.Range("K74").Value = 1400.3864
$ws.Range("L74").Value = 1635.6
$ws.Range("M74").Value = -526.3864000000001
$ws.Range("N74").Value = -3383.6

$ws.Range("H77").Value = 1460.1864
$ws.Range("I77").Value = 1400.3864
$ws.Range("J77").Value = 1635.6
$ws.Range("K77").Value = 7001.932000000001
$ws.Range("L77").Value = 8178
$ws.Range("M77").Value = -2633.932000000001
$ws.Range("N77").Value = -16914

$ws.Range("H122").Value = 1238.2307
$ws.Range("I122").Value = 955.2222
$ws.Range("J122").Value = 1875
$ws.Range("K122").Value = 2865.6666
$ws.Range("L122").Value = 5625
$ws.Range("M122").Value = -415.6666
$ws.Range("N122").Value = -10525

$ws.Range("H132").Value = 4627.13
$ws.Range("I132").Value = 5252.1274
$ws.Range("J132").Value = 3064.6365
$ws.Range("K132").Value = 15756.3822
$ws.Range("L132").Value = 9193.9095
$ws.Range("M132").Value = -13226.3822
$ws.Range("N132").Value = -14253.9095

$ws.Range("H136").Value = 2179.8262
$ws.Range("I136").Value = 1979.7778
$ws.Range("J136").Value = 2900
$ws.Range("K136").Value = 5939.3334
$ws.Range("L136").Value = 8700
$ws.Range("M136").Value = -3389.3334
$ws.Range("N136").Value = -13800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 8842.857
$ws.Range("I24").Value = 950
$ws.Range("K24").Value = 950
$ws.Range("M24").Value = -715

$ws.Range("H99").Value = 722.5
$ws.Range("I99").Value = 611
$ws.Range("J99").Value = 908.3333
$ws.Range("K99").Value = 611
$ws.Range("L99").Value = 908.3333
$ws.Range("M99").Value = 887
$ws.Range("N99").Value = -3904.3333

$ws.Range("H107").Value = 1269.3334
$ws.Range("I107").Value = 1303
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1303
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 617
$ws.Range("N107").Value = -4840

$ws.Range("H133").Value = 50000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -60120

$ws.Range("H134").Value = 6595.636
$ws.Range("I134").Value = 7724.6523
$ws.Range("J134").Value = 3998.9
$ws.Range("K134").Value = 23173.9569
$ws.Range("L134").Value = 11996.7
$ws.Range("M134").Value = -20638.9569
$ws.Range("N134").Value = -17066.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1060
$ws.Range("I16").Value = 1050
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 1050
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -763
$ws.Range("N16").Value = -1674

$ws.Range("H31").Value = 45462170
$ws.Range("I31").Value = 7884.05
$ws.Range("J31").Value = 500005000
$ws.Range("K31").Value = 7884.05
$ws.Range("L31").Value = 500005000
$ws.Range("M31").Value = -7589.05
$ws.Range("N31").Value = -500005590

$ws.Range("H34").Value = 45462170
$ws.Range("I34").Value = 7884.05
$ws.Range("J34").Value = 500005000
$ws.Range("K34").Value = 7884.05
$ws.Range("L34").Value = 500005000
$ws.Range("M34").Value = -7682.05
$ws.Range("N34").Value = -500005404

$ws.Range("H58").Value = 2071.1538
$ws.Range("I58").Value = 1905.9131
$ws.Range("K58").Value = 1905.9131
$ws.Range("M58").Value = -1702.9131

$ws.Range("H99").Value = 2178.3667
$ws.Range("I99").Value = 1885.7858
$ws.Range("J99").Value = 2434.375
$ws.Range("K99").Value = 1885.7858
$ws.Range("L99").Value = 2434.375
$ws.Range("M99").Value = -387.7858000000001
$ws.Range("N99").Value = -5430.375

$ws.Range("H113").Value = 1060
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1050
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1120
$ws.Range("N113").Value = -5440

$ws.Range("H122").Value = 968.5454999999999
$ws.Range("I122").Value = 948.5714
$ws.Range("J122").Value = 1003.5
$ws.Range("K122").Value = 2845.7142
$ws.Range("L122").Value = 3010.5
$ws.Range("M122").Value = -395.7142000000003
$ws.Range("N122").Value = -7910.5

$ws.Range("H126").Value = 2178.3667
$ws.Range("I126").Value = 1885.7858
$ws.Range("J126").Value = 2434.375
$ws.Range("K126").Value = 5657.357400000001
$ws.Range("L126").Value = 7303.125
$ws.Range("M126").Value = -3187.357400000001
$ws.Range("N126").Value = -12243.125

$ws.Range("H132").Value = 2552339.5
$ws.Range("I132").Value = 1041.9231
$ws.Range("J132").Value = 12502400
$ws.Range("K132").Value = 3125.7693
$ws.Range("L132").Value = 37507200
$ws.Range("M132").Value = -595.7692999999999
$ws.Range("N132").Value = -37512260

$ws.Range("H134").Value = 2470.61
$ws.Range("I134").Value = 2585.608
$ws.Range("J134").Value = 1737.5
$ws.Range("K134").Value = 7756.824000000001
$ws.Range("L134").Value = 5212.5
$ws.Range("M134").Value = -5221.824000000001
$ws.Range("N134").Value = -10282.5

$ws.Range("H136").Value = 2071.1538
$ws.Range("I136").Value = 1905.9131
$ws.Range("K136").Value = 5717.7393
$ws.Range("M136").Value = -3167.7393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1481.4445
$ws.Range("J34").Value = 1726.1428
$ws.Range("L34").Value = 5178.428400000001
$ws.Range("N34").Value = -5346.428400000001

$ws.Range("H39").Value = 3863.6365
$ws.Range("J39").Value = 3863.6365
$ws.Range("L39").Value = 11590.9095
$ws.Range("N39").Value = -12178.9095

$ws.Range("H55").Value = 2273.75
$ws.Range("J55").Value = 2273.75
$ws.Range("L55").Value = 6821.25
$ws.Range("N55").Value = -7175.25

$ws.Range("H131").Value = 728.95
$ws.Range("I131").Value = 309.23077
$ws.Range("J131").Value = 791.6667
$ws.Range("K131").Value = 927.69231
$ws.Range("L131").Value = 2375.0001
$ws.Range("M131").Value = 4112.30769
$ws.Range("N131").Value = -12455.0001

$ws.Range("H132").Value = 2012.6786
$ws.Range("I132").Value = 844.5
$ws.Range("J132").Value = 2479.95
$ws.Range("K132").Value = 7600.5
$ws.Range("L132").Value = 22319.55
$ws.Range("M132").Value = -5070.5
$ws.Range("N132").Value = -27379.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1874.2
$ws.Range("I102").Value = 2083.1538
$ws.Range("J102").Value = 1486.1428
$ws.Range("K102").Value = 2083.1538
$ws.Range("L102").Value = 1486.1428
$ws.Range("M102").Value = -461.1538
$ws.Range("N102").Value = -4730.1428

$ws.Range("H107").Value = 394.6154
$ws.Range("I107").Value = 531.5
$ws.Range("J107").Value = 175.6
$ws.Range("K107").Value = 531.5
$ws.Range("L107").Value = 175.6
$ws.Range("M107").Value = 1388.5
$ws.Range("N107").Value = -4015.6

$ws.Range("H113").Value = 16667862
$ws.Range("I113").Value = 50000796
$ws.Range("J113").Value = 1396.5
$ws.Range("K113").Value = 50000796
$ws.Range("L113").Value = 1396.5
$ws.Range("M113").Value = -49998626
$ws.Range("N113").Value = -5736.5

$ws.Range("H132").Value = 5506
$ws.Range("I132").Value = 5989.963
$ws.Range("J132").Value = 4054.111
$ws.Range("K132").Value = 17969.889
$ws.Range("L132").Value = 12162.333
$ws.Range("M132").Value = -15439.889
$ws.Range("N132").Value = -17222.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 57000
$ws.Range("J108").Value = 57000
$ws.Range("L108").Value = 57000
$ws.Range("N108").Value = -64680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 328.625
$ws.Range("I100").Value = 250.6923
$ws.Range("J100").Value = 666.3333
$ws.Range("K100").Value = 501.3846
$ws.Range("L100").Value = 1332.6666
$ws.Range("M100").Value = 39.61540000000002
$ws.Range("N100").Value = -2414.6666

$ws.Range("H122").Value = 4686.533
$ws.Range("I122").Value = 5380.4
$ws.Range("J122").Value = 3298.8
$ws.Range("K122").Value = 16141.2
$ws.Range("L122").Value = 9896.400000000001
$ws.Range("M122").Value = -13691.2
$ws.Range("N122").Value = -14796.4
